$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7666.3335
$ws.Range("J43").Value = 8999.5
$ws.Range("L43").Value = 8999.5
$ws.Range("N43").Value = -9137.5
$ws.Range("H49").Value = 1263.6154
$ws.Range("I49").Value = 213.57143
$ws.Range("J49").Value = 2488.6667
$ws.Range("K49").Value = 640.71429
$ws.Range("L49").Value = 7466.000100000001
$ws.Range("M49").Value = -504.71429
$ws.Range("N49").Value = -7738.000100000001
$ws.Range("H55").Value = 346.44446
$ws.Range("I55").Value = 348
$ws.Range("J55").Value = 343.33334
$ws.Range("K55").Value = 348
$ws.Range("L55").Value = 343.33334
$ws.Range("M55").Value = -134
$ws.Range("N55").Value = -771.33334
$ws.Range("H70").Value = 3040.0908
$ws.Range("I70").Value = 1999.6666
$ws.Range("J70").Value = 3430.25
$ws.Range("K70").Value = 5998.9998
$ws.Range("L70").Value = 10290.75
$ws.Range("M70").Value = -5728.9998
$ws.Range("N70").Value = -10830.75
$ws.Range("H73").Value = 3040.0908
$ws.Range("I73").Value = 1999.6666
$ws.Range("J73").Value = 3430.25
$ws.Range("K73").Value = 5998.9998
$ws.Range("L73").Value = 10290.75
$ws.Range("M73").Value = -5062.9998
$ws.Range("N73").Value = -12162.75
$ws.Range("H135").Value = 2815.2
$ws.Range("I135").Value = 3269
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 29421
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -26886
$ws.Range("N135").Value = -14070
$ws.Range("H141").Value = 2764.75
$ws.Range("I141").Value = 2727.348
$ws.Range("K141").Value = 8182.044
$ws.Range("M141").Value = -3002.044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4035096.5
$ws.Range("I32").Value = 4903435.5
$ws.Range("J32").Value = 9161.182000000001
$ws.Range("K32").Value = 4903435.5
$ws.Range("L32").Value = 9161.182000000001
$ws.Range("M32").Value = -4903148.5
$ws.Range("N32").Value = -9735.182000000001
$ws.Range("H45").Value = 1902.2354
$ws.Range("I45").Value = 1557.5454
$ws.Range("K45").Value = 1557.5454
$ws.Range("M45").Value = -1180.5454
$ws.Range("H61").Value = 1514352.4
$ws.Range("I61").Value = 5025
$ws.Range("K61").Value = 5025
$ws.Range("M61").Value = -4813
$ws.Range("H63").Value = 2342.8572
$ws.Range("I63").Value = 2600
$ws.Range("K63").Value = 2600
$ws.Range("M63").Value = -1914
$ws.Range("H66").Value = 2342.8572
$ws.Range("I66").Value = 2600
$ws.Range("K66").Value = 13000
$ws.Range("M66").Value = -9568
$ws.Range("H74").Value = 16593.281
$ws.Range("I74").Value = 964.7273
$ws.Range("K74").Value = 964.7273
$ws.Range("M74").Value = -90.72730000000001
$ws.Range("H77").Value = 16593.281
$ws.Range("I77").Value = 964.7273
$ws.Range("K77").Value = 4823.636500000001
$ws.Range("M77").Value = -455.6365000000005
$ws.Range("H136").Value = 1514352.4
$ws.Range("I136").Value = 5025
$ws.Range("K136").Value = 15075
$ws.Range("M136").Value = -12525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38960.8
$ws.Range("I20").Value = 17216.666
$ws.Range("J20").Value = 48279.715
$ws.Range("K20").Value = 17216.666
$ws.Range("L20").Value = 48279.715
$ws.Range("M20").Value = -16969.666
$ws.Range("N20").Value = -48773.715
$ws.Range("H107").Value = 1314.6666
$ws.Range("I107").Value = 1183.0588
$ws.Range("J107").Value = 1538.4
$ws.Range("K107").Value = 1183.0588
$ws.Range("L107").Value = 1538.4
$ws.Range("M107").Value = 736.9412
$ws.Range("N107").Value = -5378.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1120.6364
$ws.Range("I5").Value = 384.625
$ws.Range("K5").Value = 384.625
$ws.Range("M5").Value = -272.625
$ws.Range("H10").Value = 72099.64
$ws.Range("I10").Value = 48.666668
$ws.Range("K10").Value = 48.666668
$ws.Range("M10").Value = 90.333332
$ws.Range("H12").Value = 715644.3
$ws.Range("I12").Value = 1585
$ws.Range("J12").Value = 5000000
$ws.Range("K12").Value = 1585
$ws.Range("L12").Value = 5000000
$ws.Range("M12").Value = -1415
$ws.Range("N12").Value = -5000340
$ws.Range("H13").Value = 398.25
$ws.Range("I13").Value = 9.166667
$ws.Range("K13").Value = 9.166667
$ws.Range("M13").Value = 129.833333
$ws.Range("H86").Value = 14684.786
$ws.Range("I86").Value = 15549
$ws.Range("J86").Value = 9499.5
$ws.Range("K86").Value = 15549
$ws.Range("L86").Value = 9499.5
$ws.Range("M86").Value = -14426
$ws.Range("N86").Value = -11745.5
$ws.Range("H89").Value = 14684.786
$ws.Range("I89").Value = 15549
$ws.Range("J89").Value = 9499.5
$ws.Range("K89").Value = 77745
$ws.Range("L89").Value = 47497.5
$ws.Range("M89").Value = -72129
$ws.Range("N89").Value = -58729.5
$ws.Range("H132").Value = 36767108
$ws.Range("I132").Value = 2210.65
$ws.Range("K132").Value = 6631.950000000001
$ws.Range("M132").Value = -4101.950000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1395
$ws.Range("J68").Value = 1395
$ws.Range("L68").Value = 4185
$ws.Range("N68").Value = -5807
$ws.Range("H71").Value = 1395
$ws.Range("J71").Value = 1395
$ws.Range("L71").Value = 12555
$ws.Range("N71").Value = -20667
$ws.Range("H75").Value = 1082
$ws.Range("I75").Value = 748
$ws.Range("J75").Value = 1750
$ws.Range("K75").Value = 2244
$ws.Range("L75").Value = 5250
$ws.Range("M75").Value = -1246
$ws.Range("N75").Value = -7246
$ws.Range("H78").Value = 1082
$ws.Range("I78").Value = 748
$ws.Range("J78").Value = 1750
$ws.Range("K78").Value = 6732
$ws.Range("L78").Value = 15750
$ws.Range("M78").Value = -1740
$ws.Range("N78").Value = -25734
$ws.Range("H133").Value = 77804
$ws.Range("I133").Value = 7255
$ws.Range("J133").Value = 360000
$ws.Range("K133").Value = 21765
$ws.Range("L133").Value = 1080000
$ws.Range("M133").Value = -16705
$ws.Range("N133").Value = -1090120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 938375.6
$ws.Range("I20").Value = 3340711.8
$ws.Range("J20").Value = 37499.625
$ws.Range("K20").Value = 3340711.8
$ws.Range("L20").Value = 37499.625
$ws.Range("M20").Value = -3340466.8
$ws.Range("N20").Value = -37989.625
$ws.Range("H24").Value = 7167715.5
$ws.Range("I24").Value = 8497.5
$ws.Range("J24").Value = 8360918.5
$ws.Range("K24").Value = 8497.5
$ws.Range("L24").Value = 8360918.5
$ws.Range("M24").Value = -8324.5
$ws.Range("N24").Value = -8361264.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 10124.75
$ws.Range("I10").Value = 12999.5
$ws.Range("J10").Value = 7250
$ws.Range("K10").Value = 12999.5
$ws.Range("L10").Value = 7250
$ws.Range("M10").Value = -12859.5
$ws.Range("N10").Value = -7530
$ws.Range("H12").Value = 5500
$ws.Range("I12").Value = 4750
$ws.Range("K12").Value = 4750
$ws.Range("M12").Value = -4580
$ws.Range("H136").Value = 1005229.5
$ws.Range("I136").Value = 13238.474
$ws.Range("K136").Value = 39715.422
$ws.Range("M136").Value = -37165.422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 403.33334
$ws.Range("I5").Value = 505
$ws.Range("K5").Value = 505
$ws.Range("M5").Value = -393
$ws.Range("H6").Value = 1376.5
$ws.Range("I6").Value = 750
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -635
$ws.Range("H126").Value = 4978
$ws.Range("I126").Value = 4978
$ws.Range("K126").Value = 14934
$ws.Range("M126").Value = -12464
